# Add rule for "Non-static inner classes" to the Sheet1 code-standards table.
#
# Row 21 previously held a near-empty "Inner class" placeholder row
# (Category only, no Desc). We turn it into a fully-populated rule row:
#   - Category (A21): "Inner class" -> "Non-static inner classes"
#   - Desc     (C21): (empty) -> long description, wrapped
#   - Used     (F21): (empty) -> 0
#   - Row height grows to fit the wrapped description (34)
# Column C is widened so the new (and existing) Desc text fits better.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Desc text for the new rule (set before the Category label so the ---
# --- shared-string table gets the same ordering as the authored file) ---
$descText = "Remove  an implicit reference  of non-static inner classes (Runnable/handler/loader/task/AsyncTask) `nto their outer class(Fragment or Activity)"

$ws.Range("C21").Value = $descText
$ws.Range("C21").WrapText = $true

# --- Category label ---
$ws.Range("A21").Value = "Non-static inner classes"

# --- Used flag ---
$ws.Range("F21").Value = 0

# --- Row height to match the wrapped description ---
$ws.Rows.Item(21).RowHeight = 34

# --- Widen the Desc column (C) so the long text displays well ---
# ColumnWidth is expressed in "characters"; the host always re-derives the
# stored OOXML width using a fixed MDW-7 pixel grid (stored = round(chars*7+5)/7),
# so feed it the inverse of that so the stored width lands on 109.5.
$ws.Columns.Item(3).ColumnWidth = 109.5 - 5/7

Write-Output "Applied Non-static inner classes rule to row 21"
